# Updated cryptos list on Mon Jan 15 07:30:56 UTC 2024 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) for each coin row,
# and fixes the ordering of the WrappedEther / Chainlink rows (15 and 16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '42.651.24'
$ws.Range("E2").Value = '  -0.83%  '
# Row 3 - Ethereum
$ws.Range("D3").Value = '2.509.77'
$ws.Range("E3").Value = '  -1.65%  '
# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.46'
$ws.Range("E5").Value = '  +4.57%  '
# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.55'
$ws.Range("E6").Value = '  -2.94%  '
# Row 7 - XRP
$ws.Range("E7").Value = '  +0.67%  '
# Row 8 - USDC
$ws.Range("E8").Value = '  -0.01%  '
# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  -1.95%  '
# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.14'
$ws.Range("E10").Value = '  -1.67%  '
# Row 11 - Dogecoin
$ws.Range("E11").Value = '  -0.79%  '
# Row 12 - Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.56'
$ws.Range("E12").Value = '  -2.50%  '
# Row 13 - TRON
$ws.Range("E13").Value = '  -2.20%  '
# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '2.894.93'
$ws.Range("E14").Value = '  -1.67%  '
# Row 15 - was Chainlink, now WrappedEther (rows 15/16 swapped order)
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.554.47'
$ws.Range("E15").Value = '  +2.41%  '
# Row 16 - was WrappedEther, now Chainlink
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.54'
$ws.Range("E16").Value = '  +4.69%  '
# Row 17 - Polygon
$ws.Range("E17").Value = '  -2.80%  '
# Row 18 - WrappedBTC
$ws.Range("D18").Value = '42.637.00'
$ws.Range("E18").Value = '  -1.18%  '
# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.26'
$ws.Range("E19").Value = '  -2.51%  '
# Row 20 - ShibaInu
$ws.Range("E20").Value = '  -1.77%  '
# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.55'
$ws.Range("E21").Value = '  -1.10%  '
# Row 22 - Litecoin
$ws.Range("E22").Value = '  -0.74%  '
# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.92'
$ws.Range("E23").Value = '  -1.64%  '
# Row 24 - PancakeSwap
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.00'
$ws.Range("E24").Value = '  +0.87%  '
# Row 25 - ImmutableX
$ws.Range("E25").Value = '  -2.40%  '
# Row 26 - EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.86'
$ws.Range("E26").Value = '  -3.93%  '
# Row 27 - Dai
$ws.Range("E27").Value = '  +0.19%  '
# Row 28 - Toncoin
$ws.Range("E28").Value = '  +12.23%  '
# Row 29 - InjectiveProtocol
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.98'
$ws.Range("E29").Value = '  +2.75%  '
# Row 30 - Cosmos
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.07'
$ws.Range("E30").Value = '  -0.84%  '
# Row 31 - Filecoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.91'
$ws.Range("E31").Value = '  -2.38%  '
# Row 32 - Monero
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.93'
$ws.Range("E32").Value = '  -2.10%  '
# Row 33 - Celestia
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.48'
$ws.Range("E33").Value = '  +3.33%  '
# Row 34 - LidoDAOToken
$ws.Range("E34").Value = '  +1.20%  '
# Row 35 - ARBITRUM
$ws.Range("E35").Value = '  -3.03%  '
# Row 36 - Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0786'
$ws.Range("E36").Value = '  -2.50%  '
# Row 38 - Kaspa
$ws.Range("E38").Value = '  -1.96%  '
# Row 39 - Stellar
$ws.Range("E39").Value = '  +0.20%  '
# Row 40 - EnergySwap
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.88'
$ws.Range("E40").Value = '  -6.93%  '
# Row 41 - ApeXProtocol
$ws.Range("E41").Value = '  +2.75%  '
# Row 42 - RenderToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.87'
$ws.Range("E42").Value = '  -0.29%  '
# Row 43 - NEARProtocol
$ws.Range("E43").Value = '  -1.89%  '
# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = '  +0.10%  '
# Row 45 - VeChain
$ws.Range("E45").Value = '  -1.64%  '
# Row 46 - Maker
$ws.Range("D46").Value = '2.049.98'
$ws.Range("E46").Value = '  -2.16%  '
# Row 47 - BitcoinSV
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.44'
$ws.Range("E47").Value = '  -2.58%  '
# Row 48 - FraxShare
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.82'
$ws.Range("E48").Value = '  -2.15%  '
# Row 49 - RocketPoolETH
$ws.Range("D49").Value = '2.756.01'
$ws.Range("E49").Value = '  -1.67%  '
# Row 50 - ordi
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.03'
$ws.Range("E50").Value = '  -3.33%  '
# Row 51 - Algorand
$ws.Range("E51").Value = '  -0.76%  '
